$d = $word.ActiveDocument

# 1) Replace the word "שאל" (run text) with "תרגיל" and add rFonts hint=cs
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute("שאל", $true, $false, $false, $false, $false, $true, 1, $false, "תרגיל", 2)

# 2) Replace "ה 2" with " 2"
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute("ה 2", $true, $false, $false, $false, $false, $true, 1, $false, " 2", 2)
